$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.002699318017291303
$ws.Range("C7").Value = 1.722154581447348
$ws.Range("D7").Value = 9.984749238055725
$ws.Range("E7").Value = 3.159865382900943
$ws.Range("F7").Value = 3.202280412819196
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.05857075713622879
$ws.Range("C8").Value = 1.776051892746738
$ws.Range("D8").Value = 10.39296682505982
$ws.Range("E8").Value = 3.223812467414911
$ws.Range("F8").Value = 3.267741501677706
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = 0.1181852110595209
$ws.Range("C9").Value = 2.413718652824218
$ws.Range("D9").Value = 17.70046678524174
$ws.Range("E9").Value = 4.207192268632578
$ws.Range("F9").Value = 4.314784746857515
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -0.4634037723406197
$ws.Range("C10").Value = 2.724359023810559
$ws.Range("D10").Value = 19.81001332354748
$ws.Range("E10").Value = 4.450844113597721
$ws.Range("F10").Value = 4.607408106112586
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = -1.322354230397357
$ws.Range("C11").Value = 2.379903522002908
$ws.Range("D11").Value = 17.9112548335869
$ws.Range("E11").Value = 4.232169045960582
$ws.Range("F11").Value = 4.494807298836224
$ws.Range("G11").Value = 5
